$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-29 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-30 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("364÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "681÷3=", 2) | Out-Null
$d.Content.Find.Execute("876÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "281÷3=", 2) | Out-Null
$d.Content.Find.Execute("746÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "910÷5=", 2) | Out-Null
$d.Content.Find.Execute("511÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "766÷8=", 2) | Out-Null
$d.Content.Find.Execute("320÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "761÷7=", 2) | Out-Null
$d.Content.Find.Execute("691÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "827÷4=", 2) | Out-Null
$d.Content.Find.Execute("292÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "630÷2=", 2) | Out-Null
$d.Content.Find.Execute("512÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "353÷6=", 2) | Out-Null
$d.Content.Find.Execute("292÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "656÷9=", 2) | Out-Null
$d.Content.Find.Execute("165÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "507÷7=", 2) | Out-Null
$d.Content.Find.Execute("182÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "119÷2=", 2) | Out-Null
$d.Content.Find.Execute("273÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "230÷6=", 2) | Out-Null
$d.Content.Find.Execute("754÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "161÷6=", 2) | Out-Null
$d.Content.Find.Execute("222÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "692÷8=", 2) | Out-Null
$d.Content.Find.Execute("853÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "544÷9=", 2) | Out-Null
$d.Content.Find.Execute("960÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "594÷8=", 2) | Out-Null
$d.Content.Find.Execute("405÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "621÷7=", 2) | Out-Null
$d.Content.Find.Execute("158÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "875÷9=", 2) | Out-Null
$d.Content.Find.Execute("337÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "419÷2=", 2) | Out-Null
$d.Content.Find.Execute("681÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "872÷4=", 2) | Out-Null
$d.Content.Find.Execute("795÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "700÷3=", 2) | Out-Null
$d.Content.Find.Execute("514÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "325÷7=", 2) | Out-Null
$d.Content.Find.Execute("246÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "362÷7=", 2) | Out-Null
$d.Content.Find.Execute("678÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "374÷9=", 2) | Out-Null
$d.Content.Find.Execute("570÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "438÷3=", 2) | Out-Null
